$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.721.82'
$ws.Range("E2").Value = '  +0.22%  '
$ws.Range("D3").Value = '3.072.54'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''593.61'
$ws.Range("E5").Value = '  +0.77%  '
$ws.Range("D6").Value = '''154.78'
$ws.Range("E6").Value = '  +0.89%  '
$ws.Range("E8").Value = '  -2.75%  '
$ws.Range("D9").Value = '3.072.46'
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("E10").Value = '  -0.88%  '
$ws.Range("D11").Value = '''5.87'
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").Value = '''0.452'
$ws.Range("E12").Value = '  -2.12%  '
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").Value = '''0.0000238'
$ws.Range("E13").Value = '  -1.91%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '''36.90'
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("D15").Value = '''0.120'
$ws.Range("E15").Value = '  +1.11%  '
$ws.Range("D16").Value = '3.581.52'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").Value = '''7.21'
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("D18").Value = '63.669.57'
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").Value = '3.075.56'
$ws.Range("E19").Value = '  -0.43%  '
$ws.Range("D20").Value = '''490.39'
$ws.Range("E20").Value = '  +3.12%  '
$ws.Range("D21").Value = '''14.49'
$ws.Range("E21").Value = '  -1.79%  '
$ws.Range("D22").Value = '''0.709'
$ws.Range("E22").Value = '  -1.68%  '
$ws.Range("D23").Value = '''7.57'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '''2.47'
$ws.Range("E24").Value = '  +3.87%  '
$ws.Range("D25").Value = '''82.02'
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("D26").Value = '''12.92'
$ws.Range("E26").Value = '  -2.11%  '
$ws.Range("D27").Value = '''10.68'
$ws.Range("E27").Value = '  +8.80%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '''7.43'
$ws.Range("E29").Value = '  +2.12%  '
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").Value = '''27.34'
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("E34").Value = '  -2.03%  '
$ws.Range("D35").Value = '''1.06'
$ws.Range("E35").Value = '  +0.86%  '
$ws.Range("D36").Value = '0.0₃0823'
$ws.Range("E36").Value = '  -3.48%  '
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").Value = '''3.30'
$ws.Range("E37").Value = '  -4.13%  '
$ws.Range("B38").Value = 'Filecoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D38").Value = '''6.01'
$ws.Range("E38").Value = '  -1.99%  '
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("B40").Value = 'Cosmos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D40").Value = '''9.28'
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").Value = '''50.66'
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").Value = '''438.05'
$ws.Range("E42").Value = '  -2.04%  '
$ws.Range("D43").Value = '''0.291'
$ws.Range("E43").Value = '  +1.87%  '
$ws.Range("E44").Value = '  +4.27%  '
$ws.Range("E45").Value = '  -0.43%  '
$ws.Range("D46").Value = '2.841.04'
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").Value = '''39.42'
$ws.Range("E47").Value = '  -1.47%  '
$ws.Range("D48").Value = '''130.97'
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("D49").Value = '''25.49'
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("E51").Value = '  -1.10%  '
